# Remove the "( DenseNet is Memory hungry)" aside from the mini-batch
# bullet on slide 8, leaving just "Mini-batch size of 32".

$p   = $ppt.ActivePresentation
$s   = $p.Slides.Item(8)
$shp = $s.Shapes.Item(4)          # "TextBox 2"

# This shape auto-sizes to its text (<a:spAutoFit/>); remember the
# current height so we can restore it after the edit, since shrinking
# the text would otherwise shrink the shape too (not part of this diff).
$origHeight = $shp.Height

$tr       = $shp.TextFrame.TextRange
$fullText = $tr.Text

$oldPhrase = "Mini-batch size of 32 ( DenseNet is Memory hungry)"
$newPhrase = "Mini-batch size of 32"

$startPos = $fullText.IndexOf($oldPhrase) + 1   # 1-based for Characters()
if ($startPos -le 0) {
    throw "Could not locate target phrase in TextBox 2"
}

$target = $tr.Characters($startPos, $oldPhrase.Length)
$target.Text = $newPhrase

# Restore the shape's original height (undo the autofit re-layout).
$shp.Height = $origHeight
